$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Minimal Seats Per Party")

$ws.Range("B1").Value = "ADIK - Fighting Democratic Movement (Agonistiko Dimokratiko Kinima, ADIK)"
$ws.Range("C1").Value = "AKEL - Progressive Party of the Working People  (Anorthotiko Komma Ergazomenou Laou, AKEL)"
$ws.Range("D1").Value = "DIKO - Democratic Party  (Demokratiko Komma, DIKO)"
$ws.Range("E1").Value = "DISY - Democratic Rally  (Democratiko Synagermo, DISY)"
$ws.Range("F1").Value = "EDEK - United Central Democratic Union  (Eniaias Demokratikis Enosis Kentrou, EDEK )"
$ws.Range("G1").Value = "EDI - United Democrats  (Enomenoi Demokrates, EDI)"
$ws.Range("H1").Value = "KOP - Ecologists-Environmental Movement (Kinima Oikologon Perivallontiston, KOP)"
$ws.Range("I1").Value = "NEO - New Horizons (Neoi Orizontes, NEO)"
$ws.Range("J1").Value = "EVROKO - European Party (Evropaiko Komma, EVROKO)"
$ws.Range("K1").Value = "ELAM - National Popular Front (ELAM) (Ethiniko Laiko Metopo, ELAM)"
$ws.Range("L1").Value = "ka - Solidarity Movement (Kinima Allilengyi, ka)"
$ws.Range("M1").Value = "sypol - Citizen’s Alliance  (Symmachía Politón, sypol)"
